# Auto-generated cell updates applying the Hades_Profits diff across all 8 class sheets.
# Each block targets one worksheet; values are written directly as literals (no formulas
# in the source data), matching cell-for-cell the before/after values from the diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 5265382   # H137: 5884720.5 -> 5265382
$ws.Cells.Item(137, 9).Value = 10001740   # I137: 12501925 -> 10001740
$ws.Cells.Item(137, 11).Value = 30005220   # K137: 37505775 -> 30005220
$ws.Cells.Item(137, 13).Value = -30002670   # M137: -37503225 -> -30002670
$ws.Cells.Item(138, 8).Value = 3482762.5   # H138: 3656885 -> 3482762.5
$ws.Cells.Item(138, 9).Value = 1213.7333   # I138: 1312.8889 -> 1213.7333
$ws.Cells.Item(138, 10).Value = 6647807   # J138: 6647807.5 -> 6647807
$ws.Cells.Item(138, 11).Value = 3641.199900000001   # K138: 3938.6667 -> 3641.199900000001
$ws.Cells.Item(138, 12).Value = 19943421   # L138: 19943422.5 -> 19943421
$ws.Cells.Item(138, 13).Value = 1498.800099999999   # M138: 1201.3333 -> 1498.800099999999
$ws.Cells.Item(138, 14).Value = -19953701   # N138: -19953702.5 -> -19953701

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 1679.4   # H45: 1428.1428 -> 1679.4
$ws.Cells.Item(45, 9).Value = 1753.375   # I45: 1463.2174 -> 1753.375
$ws.Cells.Item(45, 10).Value = 1383.5   # J45: 1266.8 -> 1383.5
$ws.Cells.Item(45, 11).Value = 1753.375   # K45: 1463.2174 -> 1753.375
$ws.Cells.Item(45, 12).Value = 1383.5   # L45: 1266.8 -> 1383.5
$ws.Cells.Item(45, 13).Value = -1376.375   # M45: -1086.2174 -> -1376.375
$ws.Cells.Item(45, 14).Value = -2137.5   # N45: -2020.8 -> -2137.5
$ws.Cells.Item(61, 8).Value = 71573370   # H61: 66801852 -> 71573370
$ws.Cells.Item(61, 9).Value = 100101530   # I61: 83418160 -> 100101530
$ws.Cells.Item(61, 10).Value = 252970   # J61: 336626.66 -> 252970
$ws.Cells.Item(61, 11).Value = 100101530   # K61: 83418160 -> 100101530
$ws.Cells.Item(61, 12).Value = 252970   # L61: 336626.66 -> 252970
$ws.Cells.Item(61, 13).Value = -100101318   # M61: -83417948 -> -100101318
$ws.Cells.Item(61, 14).Value = -253394   # N61: -337050.66 -> -253394
$ws.Cells.Item(74, 8).Value = 10081241   # H74: 12601277 -> 10081241
$ws.Cells.Item(74, 9).Value = 16734441   # I74: 25101114 -> 16734441
$ws.Cells.Item(74, 11).Value = 16734441   # K74: 25101114 -> 16734441
$ws.Cells.Item(74, 13).Value = -16733567   # M74: -25100240 -> -16733567
$ws.Cells.Item(77, 8).Value = 10081241   # H77: 12601277 -> 10081241
$ws.Cells.Item(77, 9).Value = 16734441   # I77: 25101114 -> 16734441
$ws.Cells.Item(77, 11).Value = 83672205   # K77: 125505570 -> 83672205
$ws.Cells.Item(77, 13).Value = -83667837   # M77: -125501202 -> -83667837
$ws.Cells.Item(106, 8).Value = 42576.316   # H106: 49917.855 -> 42576.316
$ws.Cells.Item(106, 10).Value = 42576.316   # J106: 49917.855 -> 42576.316
$ws.Cells.Item(106, 12).Value = 42576.316   # L106: 49917.855 -> 42576.316
$ws.Cells.Item(106, 14).Value = -45100.316   # N106: -52441.855 -> -45100.316
$ws.Cells.Item(122, 8).Value = 1999.3334   # H122: 1060.1666 -> 1999.3334
$ws.Cells.Item(122, 9).Value = 1999.3334   # I122: 1060.1666 -> 1999.3334
$ws.Cells.Item(122, 11).Value = 5998.0002   # K122: 3180.4998 -> 5998.0002
$ws.Cells.Item(122, 13).Value = -3548.0002   # M122: -730.4998000000001 -> -3548.0002
$ws.Cells.Item(136, 8).Value = 71573370   # H136: 66801852 -> 71573370
$ws.Cells.Item(136, 9).Value = 100101530   # I136: 83418160 -> 100101530
$ws.Cells.Item(136, 10).Value = 252970   # J136: 336626.66 -> 252970
$ws.Cells.Item(136, 11).Value = 300304590   # K136: 250254480 -> 300304590
$ws.Cells.Item(136, 12).Value = 758910   # L136: 1009879.98 -> 758910
$ws.Cells.Item(136, 13).Value = -300302040   # M136: -250251930 -> -300302040
$ws.Cells.Item(136, 14).Value = -764010   # N136: -1014979.98 -> -764010

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 213.70589   # H22: 241.3125 -> 213.70589
$ws.Cells.Item(22, 9).Value = 180   # I22: 195.83333 -> 180
$ws.Cells.Item(22, 10).Value = 275.5   # J22: 377.75 -> 275.5
$ws.Cells.Item(22, 11).Value = 180   # K22: 195.83333 -> 180
$ws.Cells.Item(22, 12).Value = 275.5   # L22: 377.75 -> 275.5
$ws.Cells.Item(22, 13).Value = -7   # M22: -22.83332999999999 -> -7
$ws.Cells.Item(22, 14).Value = -621.5   # N22: -723.75 -> -621.5
$ws.Cells.Item(26, 8).Value = 0   # H26: 19000 -> 0
$ws.Cells.Item(26, 9).Value = 0   # I26: 19000 -> 0
$ws.Cells.Item(26, 11).Value = 0   # K26: 19000 -> 0
$ws.Cells.Item(26, 13).Value = $null   # M26: -18708 -> (blank)

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2553.282   # H31: 4114.5386 -> 2553.282
$ws.Cells.Item(31, 9).Value = 1153.96   # I31: 1530 -> 1153.96
$ws.Cells.Item(31, 10).Value = 5052.0713   # J31: 7638.909 -> 5052.0713
$ws.Cells.Item(31, 11).Value = 1153.96   # K31: 1530 -> 1153.96
$ws.Cells.Item(31, 12).Value = 5052.0713   # L31: 7638.909 -> 5052.0713
$ws.Cells.Item(31, 13).Value = -858.96   # M31: -1235 -> -858.96
$ws.Cells.Item(31, 14).Value = -5642.0713   # N31: -8228.909 -> -5642.0713
$ws.Cells.Item(34, 8).Value = 2553.282   # H34: 4114.5386 -> 2553.282
$ws.Cells.Item(34, 9).Value = 1153.96   # I34: 1530 -> 1153.96
$ws.Cells.Item(34, 10).Value = 5052.0713   # J34: 7638.909 -> 5052.0713
$ws.Cells.Item(34, 11).Value = 1153.96   # K34: 1530 -> 1153.96
$ws.Cells.Item(34, 12).Value = 5052.0713   # L34: 7638.909 -> 5052.0713
$ws.Cells.Item(34, 13).Value = -951.96   # M34: -1328 -> -951.96
$ws.Cells.Item(34, 14).Value = -5456.0713   # N34: -8042.909 -> -5456.0713
$ws.Cells.Item(125, 8).Value = 35000   # H125: 33871.23 -> 35000
$ws.Cells.Item(125, 10).Value = 35000   # J125: 33871.23 -> 35000
$ws.Cells.Item(125, 12).Value = 35000   # L125: 33871.23 -> 35000
$ws.Cells.Item(125, 14).Value = -39920   # N125: -38791.23 -> -39920
$ws.Cells.Item(132, 8).Value = 37822.035   # H132: 47987.047 -> 37822.035
$ws.Cells.Item(132, 9).Value = 1944.8096   # I132: 2421.1875 -> 1944.8096
$ws.Cells.Item(132, 10).Value = 145453.72   # J132: 169496 -> 145453.72
$ws.Cells.Item(132, 11).Value = 5834.4288   # K132: 7263.5625 -> 5834.4288
$ws.Cells.Item(132, 12).Value = 436361.16   # L132: 508488 -> 436361.16
$ws.Cells.Item(132, 13).Value = -3304.4288   # M132: -4733.5625 -> -3304.4288
$ws.Cells.Item(132, 14).Value = -441421.16   # N132: -513548 -> -441421.16

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1035.3334   # H5: 966.3158 -> 1035.3334
$ws.Cells.Item(5, 9).Value = 760   # I5: 634 -> 760
$ws.Cells.Item(5, 10).Value = 1135.4546   # J5: 1085 -> 1135.4546
$ws.Cells.Item(5, 11).Value = 2280   # K5: 1902 -> 2280
$ws.Cells.Item(5, 12).Value = 3406.3638   # L5: 3255 -> 3406.3638
$ws.Cells.Item(5, 13).Value = -2168   # M5: -1790 -> -2168
$ws.Cells.Item(5, 14).Value = -3630.3638   # N5: -3479 -> -3630.3638
$ws.Cells.Item(86, 8).Value = 1383.3334   # H86: 1549.2858 -> 1383.3334
$ws.Cells.Item(86, 10).Value = 1580   # J86: 1637.6923 -> 1580
$ws.Cells.Item(86, 12).Value = 4740   # L86: 4913.0769 -> 4740
$ws.Cells.Item(86, 14).Value = -7112   # N86: -7285.0769 -> -7112
$ws.Cells.Item(89, 8).Value = 1383.3334   # H89: 1549.2858 -> 1383.3334
$ws.Cells.Item(89, 10).Value = 1580   # J89: 1637.6923 -> 1580
$ws.Cells.Item(89, 12).Value = 14220   # L89: 14739.2307 -> 14220
$ws.Cells.Item(89, 14).Value = -26076   # N89: -26595.2307 -> -26076
$ws.Cells.Item(135, 8).Value = 1035.3334   # H135: 966.3158 -> 1035.3334
$ws.Cells.Item(135, 9).Value = 760   # I135: 634 -> 760
$ws.Cells.Item(135, 10).Value = 1135.4546   # J135: 1085 -> 1135.4546
$ws.Cells.Item(135, 11).Value = 6840   # K135: 5706 -> 6840
$ws.Cells.Item(135, 12).Value = 10219.0914   # L135: 9765 -> 10219.0914
$ws.Cells.Item(135, 13).Value = -4305   # M135: -3171 -> -4305
$ws.Cells.Item(135, 14).Value = -15289.0914   # N135: -14835 -> -15289.0914

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(46, 8).Value = 25320.666   # H46: 26025.6 -> 25320.666
$ws.Cells.Item(46, 10).Value = 25320.666   # J46: 26025.6 -> 25320.666
$ws.Cells.Item(46, 12).Value = 25320.666   # L46: 26025.6 -> 25320.666
$ws.Cells.Item(46, 14).Value = -25632.666   # N46: -26337.6 -> -25632.666
$ws.Cells.Item(57, 8).Value = 18990   # H57: 11495 -> 18990
$ws.Cells.Item(57, 9).Value = 18990   # I57: 11495 -> 18990
$ws.Cells.Item(57, 11).Value = 18990   # K57: 11495 -> 18990
$ws.Cells.Item(57, 13).Value = -18170   # M57: -10675 -> -18170
$ws.Cells.Item(100, 8).Value = 38420   # H100: 40140 -> 38420
$ws.Cells.Item(100, 10).Value = 38420   # J100: 40140 -> 38420
$ws.Cells.Item(100, 12).Value = 38420   # L100: 40140 -> 38420
$ws.Cells.Item(100, 14).Value = -40584   # N100: -42304 -> -40584
$ws.Cells.Item(101, 8).Value = 56292   # H101: 69580 -> 56292
$ws.Cells.Item(101, 10).Value = 56292   # J101: 69580 -> 56292
$ws.Cells.Item(101, 12).Value = 56292   # L101: 69580 -> 56292
$ws.Cells.Item(101, 14).Value = -62782   # N101: -76070 -> -62782
$ws.Cells.Item(102, 8).Value = 1776.6666   # H102: 2018 -> 1776.6666
$ws.Cells.Item(102, 9).Value = 1779.4736   # I102: 2019.3846 -> 1779.4736
$ws.Cells.Item(102, 10).Value = 1750   # J102: 2000 -> 1750
$ws.Cells.Item(102, 11).Value = 1779.4736   # K102: 2019.3846 -> 1779.4736
$ws.Cells.Item(102, 12).Value = 1750   # L102: 2000 -> 1750
$ws.Cells.Item(102, 13).Value = -157.4736   # M102: -397.3846000000001 -> -157.4736
$ws.Cells.Item(102, 14).Value = -4994   # N102: -5244 -> -4994
$ws.Cells.Item(122, 8).Value = 5067   # H122: 4425 -> 5067
$ws.Cells.Item(122, 9).Value = 2600   # I122: 2566.3333 -> 2600
$ws.Cells.Item(122, 11).Value = 7800   # K122: 7698.999899999999 -> 7800
$ws.Cells.Item(122, 13).Value = -5350   # M122: -5248.999899999999 -> -5350

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1872.8572   # H7: 1561.15 -> 1872.8572
$ws.Cells.Item(7, 9).Value = 1820   # I7: 1471.9412 -> 1820
$ws.Cells.Item(7, 11).Value = 1820   # K7: 1471.9412 -> 1820
$ws.Cells.Item(7, 13).Value = -1708   # M7: -1359.9412 -> -1708
$ws.Cells.Item(18, 8).Value = 200   # H18: 1333.3334 -> 200
$ws.Cells.Item(18, 9).Value = 200   # I18: 1000 -> 200
$ws.Cells.Item(18, 10).Value = 0   # J18: 2000 -> 0
$ws.Cells.Item(18, 11).Value = 200   # K18: 1000 -> 200
$ws.Cells.Item(18, 12).Value = 0   # L18: 2000 -> 0
$ws.Cells.Item(18, 13).Value = -28   # M18: -828 -> -28
$ws.Cells.Item(18, 14).Value = $null   # N18: -2344 -> (blank)
$ws.Cells.Item(20, 8).Value = 1000   # H20: 0 -> 1000
$ws.Cells.Item(20, 9).Value = 1000   # I20: 0 -> 1000
$ws.Cells.Item(20, 11).Value = 1000   # K20: 0 -> 1000
$ws.Cells.Item(20, 13).Value = -774   # M20: (none) -> -774
$ws.Cells.Item(76, 8).Value = 39664.57   # H76: 37917.176 -> 39664.57
$ws.Cells.Item(76, 10).Value = 39664.57   # J76: 37917.176 -> 39664.57
$ws.Cells.Item(76, 12).Value = 39664.57   # L76: 37917.176 -> 39664.57
$ws.Cells.Item(76, 14).Value = -40340.57   # N76: -38593.176 -> -40340.57
$ws.Cells.Item(79, 8).Value = 39664.57   # H79: 37917.176 -> 39664.57
$ws.Cells.Item(79, 10).Value = 39664.57   # J79: 37917.176 -> 39664.57
$ws.Cells.Item(79, 12).Value = 39664.57   # L79: 37917.176 -> 39664.57
$ws.Cells.Item(79, 14).Value = -42004.57   # N79: -40257.176 -> -42004.57
$ws.Cells.Item(88, 8).Value = 31750   # H88: 50555.555 -> 31750
$ws.Cells.Item(88, 9).Value = 30000   # I88: 0 -> 30000
$ws.Cells.Item(88, 10).Value = 32000   # J88: 50555.555 -> 32000
$ws.Cells.Item(88, 11).Value = 30000   # K88: 0 -> 30000
$ws.Cells.Item(88, 12).Value = 32000   # L88: 50555.555 -> 32000
$ws.Cells.Item(88, 13).Value = -29572   # M88: (none) -> -29572
$ws.Cells.Item(88, 14).Value = -32856   # N88: -51411.555 -> -32856
$ws.Cells.Item(91, 8).Value = 31750   # H91: 50555.555 -> 31750
$ws.Cells.Item(91, 9).Value = 30000   # I91: 0 -> 30000
$ws.Cells.Item(91, 10).Value = 32000   # J91: 50555.555 -> 32000
$ws.Cells.Item(91, 11).Value = 30000   # K91: 0 -> 30000
$ws.Cells.Item(91, 12).Value = 32000   # L91: 50555.555 -> 32000
$ws.Cells.Item(91, 13).Value = -28518   # M91: (none) -> -28518
$ws.Cells.Item(91, 14).Value = -34964   # N91: -53519.555 -> -34964
$ws.Cells.Item(100, 8).Value = 84831.086   # H100: 85033.336 -> 84831.086
$ws.Cells.Item(100, 9).Value = 101397.3   # I100: 168066.67 -> 101397.3
$ws.Cells.Item(100, 11).Value = 101397.3   # K100: 168066.67 -> 101397.3
$ws.Cells.Item(100, 13).Value = -100856.3   # M100: -167525.67 -> -100856.3
$ws.Cells.Item(103, 8).Value = 30806.154   # H103: 29355 -> 30806.154
$ws.Cells.Item(103, 10).Value = 30806.154   # J103: 29355 -> 30806.154
$ws.Cells.Item(103, 12).Value = 30806.154   # L103: 29355 -> 30806.154
$ws.Cells.Item(103, 14).Value = -33150.15399999999   # N103: -31699 -> -33150.15399999999
$ws.Cells.Item(122, 8).Value = 2775.4443   # H122: 2677.4 -> 2775.4443
$ws.Cells.Item(122, 9).Value = 2622.375   # I122: 2530.4443 -> 2622.375
$ws.Cells.Item(122, 11).Value = 7867.125   # K122: 7591.3329 -> 7867.125
$ws.Cells.Item(122, 13).Value = -5417.125   # M122: -5141.3329 -> -5417.125
$ws.Cells.Item(126, 8).Value = 1872.8572   # H126: 1561.15 -> 1872.8572
$ws.Cells.Item(126, 9).Value = 1820   # I126: 1471.9412 -> 1820
$ws.Cells.Item(126, 11).Value = 5460   # K126: 4415.8236 -> 5460
$ws.Cells.Item(126, 13).Value = -2990   # M126: -1945.8236 -> -2990

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(103, 8).Value = 0   # H103: 40000 -> 0
$ws.Cells.Item(103, 10).Value = 0   # J103: 40000 -> 0
$ws.Cells.Item(103, 12).Value = 0   # L103: 40000 -> 0
$ws.Cells.Item(103, 14).Value = $null   # N103: -42344 -> (blank)
$ws.Cells.Item(112, 8).Value = 34000   # H112: 31266.666 -> 34000
$ws.Cells.Item(112, 10).Value = 34000   # J112: 31266.666 -> 34000
$ws.Cells.Item(112, 12).Value = 34000   # L112: 31266.666 -> 34000
$ws.Cells.Item(112, 14).Value = -36954   # N112: -34220.666 -> -36954
$ws.Cells.Item(113, 8).Value = 832.8461   # H113: 769.5484 -> 832.8461
$ws.Cells.Item(113, 9).Value = 604.7895   # I113: 630.75 -> 604.7895
$ws.Cells.Item(113, 10).Value = 1451.8572   # J113: 917.6 -> 1451.8572
$ws.Cells.Item(113, 11).Value = 1814.3685   # K113: 1892.25 -> 1814.3685
$ws.Cells.Item(113, 12).Value = 4355.571599999999   # L113: 2752.8 -> 4355.571599999999
$ws.Cells.Item(113, 13).Value = 355.6315   # M113: 277.75 -> 355.6315
$ws.Cells.Item(113, 14).Value = -8695.571599999999   # N113: -7092.8 -> -8695.571599999999
$ws.Cells.Item(122, 8).Value = 2203.258   # H122: 2124.2122 -> 2203.258
$ws.Cells.Item(122, 9).Value = 1475.55   # I122: 1423.1364 -> 1475.55
$ws.Cells.Item(122, 11).Value = 4426.65   # K122: 4269.4092 -> 4426.65
$ws.Cells.Item(122, 13).Value = -1976.55   # M122: -1819.4092 -> -1976.55

